$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.439.33"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.644.35"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.530"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "1.876.48"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "1.652.04"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.573"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "27.417.78"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.02%  "
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "1.420.84"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.566"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").Value = "  -4.45%  "
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.824"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.17%  "
$ws.Range("D47").Value = "1.786.17"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.83%  "
